$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update)
# Set a temporary text format on the data range so numeric-looking strings
# (e.g. "0.9994", "23.511.67") are preserved verbatim instead of being
# coerced into numbers by Excel's type inference.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.511.67"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.651.15"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "300.11"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "50.60"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "0.3507"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "1.227"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "0.08064"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "0.9994"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "22.13"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "6.326"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "7.277"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "0.00001213"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "1.652.35"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "95.40"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "0.06974"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "6.641"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "12.47"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "23.523.31"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "2.420"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").Value = "3.021"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "21.13"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "151.81"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "5.185"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "131.81"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "1.836.71"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "6.911"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "2.138"
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").Value = "11.25"
$ws.Range("E34").Value = "  -7.76%  "
$ws.Range("D35").Value = "0.9915"
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("D36").Value = "0.02724"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").Value = "0.08774"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "5.950"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "0.06834"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "12.93"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "0.6912"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "1.298"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").Value = "15.68"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6405"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "0.9989"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "127.30"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07686"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "1.238"
$ws.Range("E51").Value = "  +3.00%  "

# Restore the original (unstyled) formatting so only values changed.
$dataRange.ClearFormats()
